$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Rolley"
$ws.Range("B3").Value = "Smith"
$ws.Range("B4").Value = "Balenga"
$ws.Range("B5").Value = "Issac"
$ws.Range("B6").Value = "Cruise"
$ws.Range("B7").Value = "Depp"
$ws.Range("B8").Value = "Heard"
$ws.Range("B9").Value = "Qiao"
$ws.Range("B10").Value = "Biden"
